# This script applies a row-wise permutation of the per-row "observation"
# data (Fecha, Volumen, Precio minimo/maximo/promedio, Origen, Precio $/Kg)
# for rows 2-11 of the active sheet, matching the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture current (pre-edit) values for the columns that change, keyed by row.
$D = @{}
$M = @{}
$N = @{}
$O = @{}
$P = @{}
$R = @{}
$S = @{}

for ($r = 2; $r -le 11; $r++) {
    $D[$r] = $ws.Cells.Item($r, 4).Value2   # Fecha
    $M[$r] = $ws.Cells.Item($r, 13).Value2  # Volumen
    $N[$r] = $ws.Cells.Item($r, 14).Value2  # Precio minimo
    $O[$r] = $ws.Cells.Item($r, 15).Value2  # Precio maximo
    $P[$r] = $ws.Cells.Item($r, 16).Value2  # Precio promedio ponderado
    $R[$r] = $ws.Cells.Item($r, 18).Value2  # Origen
    $S[$r] = $ws.Cells.Item($r, 19).Value2  # Precio $/Kg
}

# Row permutation: new row r gets the values that used to belong to row $perm[r]
$perm = @{
    2  = 5
    3  = 10
    4  = 8
    5  = 4
    6  = 7
    7  = 9
    8  = 11
    9  = 3
    10 = 6
    11 = 2
}

foreach ($r in 2..11) {
    $src = $perm[$r]
    $ws.Cells.Item($r, 4).Value2  = $D[$src]
    $ws.Cells.Item($r, 13).Value2 = $M[$src]
    $ws.Cells.Item($r, 14).Value2 = $N[$src]
    $ws.Cells.Item($r, 15).Value2 = $O[$src]
    $ws.Cells.Item($r, 16).Value2 = $P[$src]
    $ws.Cells.Item($r, 18).Value2 = $R[$src]
    $ws.Cells.Item($r, 19).Value2 = $S[$src]
}
